# Updated symbol list on Thu Dec 15 04:51:44 UTC 2022 with GitHub Actions
# Refreshes the latest crypto price/volume snapshot in column D (and a
# couple of "Worst in 24h" volume labels in column E) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume figures in column D are stored as text, so a leading
# apostrophe is used to stop Excel from re-interpreting them as numbers
# (which would also strip trailing/leading zeros and use scientific
# notation for the very small values).
$ws.Range("D2").Value  = "'264.94"
$ws.Range("D4").Value  = "'6.288"
$ws.Range("D6").Value  = "'3.594"
$ws.Range("D7").Value  = "'6.663"
$ws.Range("D8").Value  = "'1.345"
$ws.Range("D9").Value  = "'0.8291"
$ws.Range("D11").Value = "'0.1595"
$ws.Range("D12").Value = "'0.08211"
$ws.Range("D13").Value = "'0.03428"
$ws.Range("D15").Value = "'0.09243"
$ws.Range("D16").Value = "'3.889"
$ws.Range("D17").Value = "'0.001707"
$ws.Range("D19").Value = "'0.006242"
$ws.Range("D21").Value = "'0.001089"
$ws.Range("D23").Value = "'3.763"
$ws.Range("D26").Value = "'0.1238"
$ws.Range("D40").Value = "'0.04598"
$ws.Range("D41").Value = "'0.006959"
$ws.Range("D42").Value = "'0.1137"

$ws.Range("D43").Value = "'0.003132"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.01061"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"

$ws.Range("D45").Value = "'0.00006137"
$ws.Range("D47").Value = "'0.7787"
$ws.Range("D48").Value = "'0.1931"
